$d = $word.ActiveDocument

# 1) Update the title/heading run text (two <w:t> runs joined by a <w:br/>)
$ok1 = $d.Content.Find.Execute("Review 118: [Short] Seeing through the Brain: Image Reconstruction of Visual Perception from Human Brain Signals, 08.08.23", $false, $false, $false, $false, $false, $true, 1, $false, "Review 117b: [Short] Predicting masked tokens in stochastic locations improves masked image modeling, 07.08.23", 2)
if (-not $ok1) { throw "Could not find/replace the review title text" }
$ok2 = $d.Content.Find.Execute("https://huggingface.co/papers/2308.02510", $false, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2308.00566.pdf", 2)
if (-not $ok2) { throw "Could not find/replace the huggingface paper URL" }

# 2) Update the bold "Paper:" line
$ok3 = $d.Content.Find.Execute("Paper: https://arxiv.org/abs/2208.03666v4", $false, $false, $false, $false, $false, $true, 1, $false, "Paper: https://arxiv.org/abs/2303.00289v1", 2)
if (-not $ok3) { throw "Could not find/replace the 'Paper:' arxiv link" }

# 3) Remove the extra empty "Normal"-styled paragraph (4th paragraph)
$d.Paragraphs(4).Range.Delete()

# 4) Replace the long Hebrew review paragraph body (now paragraph 4) with new content,
#    preserving its paragraph style/properties by targeting a range that excludes the
#    paragraph mark, and inserting fresh run XML (text + <w:br/> line breaks) via InsertXML.
$body = $d.Paragraphs(4)
$bodyRange = $d.Range($body.Range.Start, $body.Range.End - 1)
$xml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:r><w:t xml:space='preserve'>היום ב-#shorthebrewpapereviews סוקרים מאמר של כמה חוקרים ישראלים עם Yann LeCun האגדי!! שיטות למידה self-supervised (או SSL) הפכו להיות מאוד פופולריות לבניית ייצוג עוצמתי עבור דאטה ויזואלי (תמונות) שניתן להשתמש בו למשימות מגוונות. שיטות אלו לא דורשות דאטה מתויג ולכן ניתן לאמן אותם על דאטהסטים ענקיים של תמונות מהאינטרנט. </w:t><w:br/><w:br/><w:t xml:space='preserve'>בדרך כלל שיטת SSL מהנדסת משימה שלא דורשת תמונות מתויגות. למשל אחד המאמר האחרונים של יאן לקון (I-JEPA) המשימה הייתה חיזוי הייצוג (embedding) של פאץ בתמונה נתונה בהינתן ייצוגים של פאצ'ים אחרים של התמונה. ככה ייצוג שנבנה לומד להפיק את המאפיינים הסמנטיים של הפאצ'ים מייצוג הויזואלי של הפאצ'ים באותה תמונה. במאמר I-JEPA המודל מקבל את הייצוגים של כמה פאצ'ים (ההקשר) יחד עם הקידוד המיקום שלו בתמונה (positional encoding) של המיקומים של הפאץ' שחיזויו היה צריך לחזות היה מיוצג עם וקטור המיסוך (הקבוע עבור כל הפאצ'ים) וגם קידוד המיקום שלו בתמונה. </w:t><w:br/><w:br/><w:t>במאמר הנסקר המחברים מבקשים להכליל את הגישה של I-JEPA ובמקום קידוד מיקום מדויק להעביר למודל קידוק מקומי מורעש (גם עבור פאצ'י הקשר וגם עבור פאצ'י שייצוגם נחזים). איך זה נעשה? פשוט מוסיפים וקטור גאוסי עם מטריצת קווריאנס S נלמדת לוקטור קידוד מיקום. ככה אנו הופכים את משימת SSL מורכבת יותר וכתוצאה מכך הייצוגים המופקים באמצעותה משתפרים. מאחר וצריך ללמוד פרמטרים של ההתפלגות שממנה צריך לדגום את הוקטור המורעש המייצג מיקום אז בדומה ל-VAE משתמשים ב-reparameterization trick.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$bodyRange.InsertXML($xml)

if ($d.Paragraphs.Count -ne 5) { throw "Unexpected paragraph count after edit: $($d.Paragraphs.Count)" }

Write-Output "Done. Paragraph count:"
Write-Output $d.Paragraphs.Count
